# Add season record columns (Wins, Losses, Ties) to the DET_2007 sheet.
#
# The sheet currently spans A1:AC46 (a header row plus 45 data/header-repeat
# rows). We append three new columns - AD (Wins), AE (Losses), AF (Ties) -
# with a header label in row 1 and the team's season record (88-74-0) repeated
# on every remaining row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 46

# Copy the header formatting (bold font, border, centered alignment) from the
# existing last header cell (AC1) onto the three new header cells so they
# reuse the same style index as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every remaining row (data rows + the repeated
# header row at r=46), matching the team's 88-74-0 record.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
